$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the dSF column (F) for rows 3-6 per repulled data / mean calculation
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = -4
